# Append two new diary-entry paragraphs after the existing content,
# matching the "2022年6月X日星期X" / weather-note paragraph pattern
# already used throughout the document.

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Paragraph 1: date line, split into three runs just like the other
# date paragraphs in the document ("2" / "022" / "年6月7日星期二").
$dateParagraph = '<w:p ' + $wNs + '>' `
  + '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>2</w:t></w:r>' `
  + '<w:r><w:t>022</w:t></w:r>' `
  + '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>年6月7日星期二</w:t></w:r>' `
  + '</w:p>'

# Paragraph 2: weather / diary note for that day.
$noteParagraph = '<w:p ' + $wNs + '>' `
  + '<w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>' `
  + '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>晴，今天是高考的第一天，上午考语文，下午考数学</w:t></w:r>' `
  + '</w:p>'

# Collapse to the end of the document body and insert the new paragraphs
# right before the final section break, after the last existing paragraph.
$endOfDoc = $d.Content
$endOfDoc.Collapse(0)
$endOfDoc.InsertXML($dateParagraph + $noteParagraph)
